# Update cryptocurrency price/volume data (D and E columns) to match latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect column D from Excel's automatic text-to-number coercion
# (values like "596.79" would otherwise be parsed as numbers, while
# values like "67.358.78" survive as text already). Force text format,
# write the values, then clear the temporary formatting so the cells
# end up unstyled again, matching the original workbook.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "67.358.78"
$ws.Range("D3").Value = "3.488.55"
$ws.Range("D5").Value = "596.79"
$ws.Range("D6").Value = "179.70"
$ws.Range("D8").Value = "0.606"
$ws.Range("D9").Value = "3.490.34"
$ws.Range("D13").Value = "4.092.24"
$ws.Range("D14").Value = "32.23"
$ws.Range("D16").Value = "67.350.14"
$ws.Range("D18").Value = "3.490.28"
$ws.Range("D20").Value = "14.27"
$ws.Range("D21").Value = "389.48"
$ws.Range("D22").Value = "7.92"
$ws.Range("D23").Value = "74.03"
$ws.Range("D24").Value = "0.541"
$ws.Range("D25").Value = "0.999"
$ws.Range("D27").Value = "0.0000121"
$ws.Range("D28").Value = "10.36"
$ws.Range("D30").Value = "0.979"
$ws.Range("D31").Value = "6.19"
$ws.Range("D32").Value = "1.42"
$ws.Range("D34").Value = "23.53"
$ws.Range("D35").Value = "7.40"
$ws.Range("D38").Value = "163.49"
$ws.Range("D41").Value = "1.88"
$ws.Range("D42").Value = "6.82"
$ws.Range("D44").Value = "2.848.67"
$ws.Range("D45").Value = "26.32"
$ws.Range("D46").Value = "26.72"
$ws.Range("D47").Value = "0.0721"
$ws.Range("D48").Value = "41.69"
$ws.Range("D50").Value = "333.32"

$dRange.ClearFormats()

# Column E values are percentage strings padded with spaces, already
# safe as plain text.
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("E6").Value = "  +4.06%  "
$ws.Range("E8").Value = "  +2.83%  "
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("E10").Value = "  +4.73%  "
$ws.Range("E11").Value = "  -1.78%  "
$ws.Range("E12").Value = "  +1.27%  "
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("E14").Value = "  +9.87%  "
$ws.Range("E15").Value = "  +0.54%  "
$ws.Range("E16").Value = "  +0.73%  "
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("E20").Value = "  +0.32%  "
$ws.Range("E21").Value = "  -1.32%  "
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("E23").Value = "  +1.06%  "
$ws.Range("E24").Value = "  +1.32%  "
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("E26").Value = "  +0.93%  "
$ws.Range("E27").Value = "  +0.19%  "
$ws.Range("E28").Value = "  +1.50%  "
$ws.Range("E29").Value = "  -2.78%  "
$ws.Range("E30").Value = "  -1.86%  "
$ws.Range("E31").Value = "  +0.72%  "
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("E33").Value = "  +0.71%  "
$ws.Range("E34").Value = "  -0.72%  "
$ws.Range("E35").Value = "  +0.54%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  -0.44%  "
$ws.Range("E38").Value = "  +0.65%  "
$ws.Range("E39").Value = "  -0.82%  "
$ws.Range("E40").Value = "  +10.48%  "
$ws.Range("E41").Value = "  -1.02%  "
$ws.Range("E42").Value = "  -0.90%  "
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("E44").Value = "  +0.52%  "
$ws.Range("E45").Value = "  +0.75%  "
$ws.Range("E46").Value = "  -1.53%  "
$ws.Range("E47").Value = "  -2.15%  "
$ws.Range("E48").Value = "  -2.33%  "
$ws.Range("E49").Value = "  -0.57%  "
$ws.Range("E50").Value = "  -1.18%  "
$ws.Range("E51").Value = "  -1.56%  "
